$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Rows.Item(24).Insert() | Out-Null
$ws.Cells.Item(23, 1).Copy() | Out-Null
$ws.Cells.Item(24, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# row 2: 景德镇·BM次元盛典运动番only（取消）
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = '''2024-06-22'
$ws.Cells.Item(2, 3).Value = '''景德镇·BM次元盛典运动番only（取消）'
$ws.Cells.Item(2, 4).Value = '''广场南路金幕影城旁 罗曼园宴会酒店'
$ws.Cells.Item(2, 5).Value = '''2024.06.22 10:00-06.22 17:00'
$ws.Cells.Item(2, 6).Value = 204
$ws.Cells.Item(2, 7).Value = '不可售'
$ws.Cells.Item(2, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85197'
$ws.Cells.Item(2, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202404/Z6eXz0su1714292081978.png'

# row 3: 上饶·BM次元盛典运动番only（取消）
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = '''2024-06-23'
$ws.Cells.Item(3, 3).Value = '''上饶·BM次元盛典运动番only（取消）'
$ws.Cells.Item(3, 4).Value = '''春江北大道时光PARK内 博悦宴会艺术中心'
$ws.Cells.Item(3, 5).Value = '''2024.06.23 10:00-06.23 17:00'
$ws.Cells.Item(3, 6).Value = 279
$ws.Cells.Item(3, 7).Value = '不可售'
$ws.Cells.Item(3, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85201'
$ws.Cells.Item(3, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202404/30dgkbjT1714293499693.png'

# row 4: 赣州·清风霁月·光夜only
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = '''2024-06-23'
$ws.Cells.Item(4, 3).Value = '''赣州·清风霁月·光夜only'
$ws.Cells.Item(4, 4).Value = '''平安大道 麋鹿LiveHouse'
$ws.Cells.Item(4, 5).Value = '''2024.06.23 14:00-06.23 20:00'
$ws.Cells.Item(4, 6).Value = 84
$ws.Cells.Item(4, 7).Value = 158
$ws.Cells.Item(4, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86993'
$ws.Cells.Item(4, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/PklWR8EP1717429316070.jpeg'

# row 5: 南昌·第五人格only
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = '''2024-06-29'
$ws.Cells.Item(5, 3).Value = '''南昌·第五人格only'
$ws.Cells.Item(5, 4).Value = '''高处见美好生活公园 百家喜宴高新店'
$ws.Cells.Item(5, 5).Value = '''2024.06.29 10:00-06.29 17:00'
$ws.Cells.Item(5, 6).Value = 293
$ws.Cells.Item(5, 7).Value = 65
$ws.Cells.Item(5, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87043'
$ws.Cells.Item(5, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202405/zir2PYz81717071721569.jpeg'

# row 6: 萍乡·BM次元盛典运动番only
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = '''2024-06-29'
$ws.Cells.Item(6, 3).Value = '''萍乡·BM次元盛典运动番only'
$ws.Cells.Item(6, 4).Value = '''康庄路3号 萍乡梅园国际大酒店'
$ws.Cells.Item(6, 5).Value = '''2024.06.29 10:00-06.29 17:00'
$ws.Cells.Item(6, 6).Value = 272
$ws.Cells.Item(6, 7).Value = 55
$ws.Cells.Item(6, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85192'
$ws.Cells.Item(6, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png'

# row 7: 南昌·ChinastyleCOSPLAY  
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = '''2024-06-30'
$ws.Cells.Item(7, 3).Value = '''南昌·ChinastyleCOSPLAY  '
$ws.Cells.Item(7, 4).Value = '''真君路888号 南昌华侨城玩美公园'
$ws.Cells.Item(7, 5).Value = '''2024.06.30 09:30-07.02 17:30'
$ws.Cells.Item(7, 6).Value = 116
$ws.Cells.Item(7, 7).Value = 65
$ws.Cells.Item(7, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87045'
$ws.Cells.Item(7, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/wajWy7ID1717149642528.jpeg'

# row 8: 宜春·BM次元盛典运动番only
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = '''2024-06-30'
$ws.Cells.Item(8, 3).Value = '''宜春·BM次元盛典运动番only'
$ws.Cells.Item(8, 4).Value = '''鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)'
$ws.Cells.Item(8, 5).Value = '''2024.06.30 10:00-06.30 17:00'
$ws.Cells.Item(8, 6).Value = 265
$ws.Cells.Item(8, 7).Value = 55
$ws.Cells.Item(8, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=84636'
$ws.Cells.Item(8, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png'

# row 9: 南昌·次元星球动漫游戏展
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = '''2024-07-06'
$ws.Cells.Item(9, 3).Value = '''南昌·次元星球动漫游戏展'
$ws.Cells.Item(9, 4).Value = '''龙蟠街666号融创茂1层 融创茂'
$ws.Cells.Item(9, 5).Value = '''2024.07.06 10:00-07.06 17:00'
$ws.Cells.Item(9, 6).Value = 22
$ws.Cells.Item(9, 7).Value = '不可售'
$ws.Cells.Item(9, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86405'
$ws.Cells.Item(9, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg'

# row 10: 鹰潭·BM次元盛典运动番only
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = '''2024-07-06'
$ws.Cells.Item(10, 3).Value = '''鹰潭·BM次元盛典运动番only'
$ws.Cells.Item(10, 4).Value = '''体育馆东路2号九小隔壁 忆江南•宴会楼'
$ws.Cells.Item(10, 5).Value = '''2024.07.06 10:00-07.06 17:00'
$ws.Cells.Item(10, 6).Value = 52
$ws.Cells.Item(10, 7).Value = 55
$ws.Cells.Item(10, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85997'
$ws.Cells.Item(10, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png'

# row 11: 赣州·BM次元盛典运动番only
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = '''2024-07-07'
$ws.Cells.Item(11, 3).Value = '''赣州·BM次元盛典运动番only'
$ws.Cells.Item(11, 4).Value = '''米瑞金路2口0号上客天下1楼 上客天下.老虔州'
$ws.Cells.Item(11, 5).Value = '''2024.07.07 10:00-07.07 17:00'
$ws.Cells.Item(11, 6).Value = 38
$ws.Cells.Item(11, 7).Value = 55
$ws.Cells.Item(11, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86602'
$ws.Cells.Item(11, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png'

# row 12: 新余·2024第三届MG动漫嘉年华
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = '''2024-07-12'
$ws.Cells.Item(12, 3).Value = '''新余·2024第三届MG动漫嘉年华'
$ws.Cells.Item(12, 4).Value = '''仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅'
$ws.Cells.Item(12, 5).Value = '''2024.07.12 10:00-07.13 17:30'
$ws.Cells.Item(12, 6).Value = 123
$ws.Cells.Item(12, 7).Value = 55
$ws.Cells.Item(12, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86536'
$ws.Cells.Item(12, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg'

# row 13: 南昌·SuperComic动漫游戏博览会
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = '''2024-07-13'
$ws.Cells.Item(13, 3).Value = '''南昌·SuperComic动漫游戏博览会'
$ws.Cells.Item(13, 4).Value = '''怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Cells.Item(13, 5).Value = '''2024.07.13 09:00-07.14 17:00'
$ws.Cells.Item(13, 6).Value = 2547
$ws.Cells.Item(13, 7).Value = 65
$ws.Cells.Item(13, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86992'
$ws.Cells.Item(13, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/wQTAjelJ1717642148929.jpeg'

# row 14: 南昌·SuperComic配音演员刘明月专场见面会
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = '''2024-07-13'
$ws.Cells.Item(14, 3).Value = '''南昌·SuperComic配音演员刘明月专场见面会'
$ws.Cells.Item(14, 4).Value = '''怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Cells.Item(14, 5).Value = '''2024.07.13 09:00-07.13 17:00'
$ws.Cells.Item(14, 6).Value = 88
$ws.Cells.Item(14, 7).Value = 168
$ws.Cells.Item(14, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87570'
$ws.Cells.Item(14, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202406/1D1reIl81718609013880.png'

# row 15: 南昌·THO-梦违赣鄱荟萃·叁~幻想Strawberry~!!
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = '''2024-07-13'
$ws.Cells.Item(15, 3).Value = '''南昌·THO-梦违赣鄱荟萃·叁~幻想Strawberry~!!'
$ws.Cells.Item(15, 4).Value = '''民德路411号 东方豪景花园酒店(民德路店)'
$ws.Cells.Item(15, 5).Value = '''2024.07.13 09:30-07.13 17:30'
$ws.Cells.Item(15, 6).Value = 30
$ws.Cells.Item(15, 7).Value = 65
$ws.Cells.Item(15, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87668'
$ws.Cells.Item(15, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/Bk9cYryT1718360290362.jpeg'

# row 16: 宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = '''2024-07-13'
$ws.Cells.Item(16, 3).Value = '''宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华'
$ws.Cells.Item(16, 4).Value = '''宜春国际商贸城会展中心 宜春国际商贸城会展中心'
$ws.Cells.Item(16, 5).Value = '''2024.07.13 10:00-07.14 17:00'
$ws.Cells.Item(16, 6).Value = 61
$ws.Cells.Item(16, 7).Value = 55
$ws.Cells.Item(16, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86667'
$ws.Cells.Item(16, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg'

# row 17: 赣州·十万伏特-次元交流会（夏）
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = '''2024-07-13'
$ws.Cells.Item(17, 3).Value = '''赣州·十万伏特-次元交流会（夏）'
$ws.Cells.Item(17, 4).Value = '''梅关大道36-16号 麋鹿星球艺术中心'
$ws.Cells.Item(17, 5).Value = '''2024.07.13 09:30-07.13 17:00'
$ws.Cells.Item(17, 6).Value = 19
$ws.Cells.Item(17, 7).Value = 45
$ws.Cells.Item(17, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87597'
$ws.Cells.Item(17, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/87yQ4Hmf1718681348727.jpeg'

# row 18: 南昌·赛马娘ONLY
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = '''2024-07-14'
$ws.Cells.Item(18, 3).Value = '''南昌·赛马娘ONLY'
$ws.Cells.Item(18, 4).Value = '''洪城路99号 锦都皇冠酒店(八一广场火车站店)'
$ws.Cells.Item(18, 5).Value = '''2024.07.14 09:00-07.14 17:30'
$ws.Cells.Item(18, 6).Value = 42
$ws.Cells.Item(18, 7).Value = 68
$ws.Cells.Item(18, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87367'
$ws.Cells.Item(18, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/wXQuIKtu1718165450704.png'

# row 19: 吉安·COMIC LIFE次元假日05
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = '''2024-07-14'
$ws.Cells.Item(19, 3).Value = '''吉安·COMIC LIFE次元假日05'
$ws.Cells.Item(19, 4).Value = '''东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws.Cells.Item(19, 5).Value = '''2024.07.14 09:00-07.14 18:00'
$ws.Cells.Item(19, 6).Value = 542
$ws.Cells.Item(19, 7).Value = 52.1
$ws.Cells.Item(19, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85924'
$ws.Cells.Item(19, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg'

# row 20: 赣州·第四届赣州半夏动漫展
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = '''2024-07-19'
$ws.Cells.Item(20, 3).Value = '''赣州·第四届赣州半夏动漫展'
$ws.Cells.Item(20, 4).Value = '''105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws.Cells.Item(20, 5).Value = '''2024.07.19 10:00-07.21 17:00'
$ws.Cells.Item(20, 6).Value = 605
$ws.Cells.Item(20, 7).Value = 55
$ws.Cells.Item(20, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86587'
$ws.Cells.Item(20, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg'

# row 21: 南昌·漫拥动漫嘉年华Pro-追光启航
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = '''2024-07-20'
$ws.Cells.Item(21, 3).Value = '''南昌·漫拥动漫嘉年华Pro-追光启航'
$ws.Cells.Item(21, 4).Value = '''小蓝南路420号 洪州体育馆'
$ws.Cells.Item(21, 5).Value = '''2024.07.20 09:00-07.21 17:00'
$ws.Cells.Item(21, 6).Value = 181
$ws.Cells.Item(21, 7).Value = 52.5
$ws.Cells.Item(21, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85796'
$ws.Cells.Item(21, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png'

# row 22: 乐平·CY境界次元动漫夏时庆
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = '''2024-07-21'
$ws.Cells.Item(22, 3).Value = '''乐平·CY境界次元动漫夏时庆'
$ws.Cells.Item(22, 4).Value = '''翥山西路182号 佳佳基大酒店'
$ws.Cells.Item(22, 5).Value = '''2024.07.21 10:00-07.21 17:00'
$ws.Cells.Item(22, 6).Value = 93
$ws.Cells.Item(22, 7).Value = 35
$ws.Cells.Item(22, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86768'
$ws.Cells.Item(22, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png'

# row 23: 九江·SXD动漫嘉年华
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = '''2024-07-21'
$ws.Cells.Item(23, 3).Value = '''九江·SXD动漫嘉年华'
$ws.Cells.Item(23, 4).Value = '''湓浦街道大中路339号 百嘉洲际酒店'
$ws.Cells.Item(23, 5).Value = '''2024.07.21 10:00-07.21 17:30'
$ws.Cells.Item(23, 6).Value = 51
$ws.Cells.Item(23, 7).Value = 45
$ws.Cells.Item(23, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86832'
$ws.Cells.Item(23, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg'

# row 24: 抚州·临次元08·盛夏动漫狂欢节
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = '''2024-07-21'
$ws.Cells.Item(24, 3).Value = '''抚州·临次元08·盛夏动漫狂欢节'
$ws.Cells.Item(24, 4).Value = '''伍塘路1098号 乐课篮球公园'
$ws.Cells.Item(24, 5).Value = '''2024.07.21 10:00-07.21 16:00'
$ws.Cells.Item(24, 6).Value = 11
$ws.Cells.Item(24, 7).Value = 39.9
$ws.Cells.Item(24, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87763'
$ws.Cells.Item(24, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/6qgetbCh1718720523395.jpeg'

# row 25: 萍乡·NL14动漫游戏展·夏日狂想曲
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = '''2024-07-21'
$ws.Cells.Item(25, 3).Value = '''萍乡·NL14动漫游戏展·夏日狂想曲'
$ws.Cells.Item(25, 4).Value = '''公园南路168号(近工行城北分理处) 梅生嘉华酒店'
$ws.Cells.Item(25, 5).Value = '''2024.07.21 10:00-07.21 17:00'
$ws.Cells.Item(25, 6).Value = 55
$ws.Cells.Item(25, 7).Value = 40
$ws.Cells.Item(25, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86658'
$ws.Cells.Item(25, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg'

# row 26: 南昌·萌卡动漫展
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = '''2024-07-26'
$ws.Cells.Item(26, 3).Value = '''南昌·萌卡动漫展'
$ws.Cells.Item(26, 4).Value = '''八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$ws.Cells.Item(26, 5).Value = '''2024.07.26 09:00-07.28 17:00'
$ws.Cells.Item(26, 6).Value = 2141
$ws.Cells.Item(26, 7).Value = 39.9
$ws.Cells.Item(26, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86776'
$ws.Cells.Item(26, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg'

# row 27: 江西·次元星河动漫游戏嘉年华
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = '''2024-07-27'
$ws.Cells.Item(27, 3).Value = '''江西·次元星河动漫游戏嘉年华'
$ws.Cells.Item(27, 4).Value = '''九龙大道1177号 南昌绿地国际博览中心'
$ws.Cells.Item(27, 5).Value = '''2024.07.27 10:00-07.28 17:00'
$ws.Cells.Item(27, 6).Value = 4240
$ws.Cells.Item(27, 7).Value = 69
$ws.Cells.Item(27, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85493'
$ws.Cells.Item(27, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png'

# row 28: 赣州·马娘only
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = '''2024-07-27'
$ws.Cells.Item(28, 3).Value = '''赣州·马娘only'
$ws.Cells.Item(28, 4).Value = '''火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)'
$ws.Cells.Item(28, 5).Value = '''2024.07.27 09:00-07.27 17:00'
$ws.Cells.Item(28, 6).Value = 32
$ws.Cells.Item(28, 7).Value = 60
$ws.Cells.Item(28, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86772'
$ws.Cells.Item(28, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png'

# row 29: 赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = '''2024-07-28'
$ws.Cells.Item(29, 3).Value = '''赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会'
$ws.Cells.Item(29, 4).Value = '''兴国路恒大帝景西门 江西长庚控股有限公司'
$ws.Cells.Item(29, 5).Value = '''2024.07.28 11:00-07.28 17:00'
$ws.Cells.Item(29, 6).Value = 66
$ws.Cells.Item(29, 7).Value = 56
$ws.Cells.Item(29, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85688'
$ws.Cells.Item(29, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png'

# row 30: 宜春·第三十五届静卿国风动漫文化展览会
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = '''2024-07-30'
$ws.Cells.Item(30, 3).Value = '''宜春·第三十五届静卿国风动漫文化展览会'
$ws.Cells.Item(30, 4).Value = '''宜阳大道19号(交通银行旁) 宜春安缦文华酒店'
$ws.Cells.Item(30, 5).Value = '''2024.07.30 09:00-07.30 17:00'
$ws.Cells.Item(30, 6).Value = 467
$ws.Cells.Item(30, 7).Value = 45
$ws.Cells.Item(30, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86684'
$ws.Cells.Item(30, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/45bGPXfQ1716709212619.jpeg'

# row 31: 南昌·幻梦境国际动漫游戏嘉年华1th
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = '''2024-08-03'
$ws.Cells.Item(31, 3).Value = '''南昌·幻梦境国际动漫游戏嘉年华1th'
$ws.Cells.Item(31, 4).Value = '''南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Cells.Item(31, 5).Value = '''2024.08.03 09:00-08.04 17:30'
$ws.Cells.Item(31, 6).Value = 1234
$ws.Cells.Item(31, 7).Value = 64
$ws.Cells.Item(31, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=83980'
$ws.Cells.Item(31, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg'

# row 32: 吉安·COMIC LIFE周年庆典
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = '''2024-08-03'
$ws.Cells.Item(32, 3).Value = '''吉安·COMIC LIFE周年庆典'
$ws.Cells.Item(32, 4).Value = '''东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws.Cells.Item(32, 5).Value = '''2024.08.03 09:30-08.03 18:00'
$ws.Cells.Item(32, 6).Value = 246
$ws.Cells.Item(32, 7).Value = 46.6
$ws.Cells.Item(32, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87164'
$ws.Cells.Item(32, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/NWD9iQ9h1717598526259.jpeg'

# row 33: 景德镇·第十五届瓷都ACG动漫游戏博览会
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = '''2024-08-03'
$ws.Cells.Item(33, 3).Value = '''景德镇·第十五届瓷都ACG动漫游戏博览会'
$ws.Cells.Item(33, 4).Value = '''迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Cells.Item(33, 5).Value = '''2024.08.03 09:00-08.04 17:00'
$ws.Cells.Item(33, 6).Value = 2140
$ws.Cells.Item(33, 7).Value = 55
$ws.Cells.Item(33, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86341'
$ws.Cells.Item(33, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png'

# row 34: 景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = '''2024-08-03'
$ws.Cells.Item(34, 3).Value = '''景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票'
$ws.Cells.Item(34, 4).Value = '''迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Cells.Item(34, 5).Value = '''2024.08.03 08:30-08.03 17:00'
$ws.Cells.Item(34, 6).Value = 568
$ws.Cells.Item(34, 7).Value = '已售罄'
$ws.Cells.Item(34, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85981'
$ws.Cells.Item(34, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png'

# row 35: 樟树·第二届静卿国风动漫文化展览会
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = '''2024-08-03'
$ws.Cells.Item(35, 3).Value = '''樟树·第二届静卿国风动漫文化展览会'
$ws.Cells.Item(35, 4).Value = '''杏佛路89号 樟树银河国际酒店'
$ws.Cells.Item(35, 5).Value = '''2024.08.03 09:00-08.03 17:00'
$ws.Cells.Item(35, 6).Value = 474
$ws.Cells.Item(35, 7).Value = 45
$ws.Cells.Item(35, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86683'
$ws.Cells.Item(35, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg'

# row 36: 萍乡·AU9夏至国漫展
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = '''2024-08-03'
$ws.Cells.Item(36, 3).Value = '''萍乡·AU9夏至国漫展'
$ws.Cells.Item(36, 4).Value = '''金陵东路18号 萍乡市体育馆'
$ws.Cells.Item(36, 5).Value = '''2024.08.03 10:00-08.03 17:00'
$ws.Cells.Item(36, 6).Value = 69
$ws.Cells.Item(36, 7).Value = 45
$ws.Cells.Item(36, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86453'
$ws.Cells.Item(36, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg'

# row 37: 赣州·第一届环梦动漫游戏嘉年华
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = '''2024-08-03'
$ws.Cells.Item(37, 3).Value = '''赣州·第一届环梦动漫游戏嘉年华'
$ws.Cells.Item(37, 4).Value = '''105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws.Cells.Item(37, 5).Value = '''2024.08.03 09:00-08.05 17:00'
$ws.Cells.Item(37, 6).Value = 25
$ws.Cells.Item(37, 7).Value = 36.6
$ws.Cells.Item(37, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87449'
$ws.Cells.Item(37, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/iC3PuUfR1717740188790.jpeg'

# row 38: 上饶·第十五届IX Group国风嘉年华暨十周年庆典
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = '''2024-08-04'
$ws.Cells.Item(38, 3).Value = '''上饶·第十五届IX Group国风嘉年华暨十周年庆典'
$ws.Cells.Item(38, 4).Value = '''高铁经济试验区凤凰东大道16号7幢 上饶饶商金茂诚悦酒店(上饶高铁站)'
$ws.Cells.Item(38, 5).Value = '''2024.08.04 09:30-08.04 17:30'
$ws.Cells.Item(38, 6).Value = 135
$ws.Cells.Item(38, 7).Value = 60
$ws.Cells.Item(38, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87225'
$ws.Cells.Item(38, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202406/l5fIXZSX1717562269098.jpeg'

# row 39: 九江·第一届异次元动漫嘉年华
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = '''2024-08-04'
$ws.Cells.Item(39, 3).Value = '''九江·第一届异次元动漫嘉年华'
$ws.Cells.Item(39, 4).Value = '''长虹西大道兴城广场99号 九江半岛宾馆'
$ws.Cells.Item(39, 5).Value = '''2024.08.04 08:00-08.04 17:00'
$ws.Cells.Item(39, 6).Value = 296
$ws.Cells.Item(39, 7).Value = 45
$ws.Cells.Item(39, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=84407'
$ws.Cells.Item(39, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202406/65hJjOfJ1717642614493.jpeg'

# row 40: 南昌·第一届异次元动漫嘉年华
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = '''2024-08-06'
$ws.Cells.Item(40, 3).Value = '''南昌·第一届异次元动漫嘉年华'
$ws.Cells.Item(40, 4).Value = '''民德路411号 东方豪景花园酒店(民德路店)'
$ws.Cells.Item(40, 5).Value = '''2024.08.06 08:00-08.06 17:00'
$ws.Cells.Item(40, 6).Value = 442
$ws.Cells.Item(40, 7).Value = 55
$ws.Cells.Item(40, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=84102'
$ws.Cells.Item(40, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg'

# row 41: 赣州·第二届异次元动漫嘉年华
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = '''2024-08-08'
$ws.Cells.Item(41, 3).Value = '''赣州·第二届异次元动漫嘉年华'
$ws.Cells.Item(41, 4).Value = '''金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆'
$ws.Cells.Item(41, 5).Value = '''2024.08.08 08:00-08.08 17:00'
$ws.Cells.Item(41, 6).Value = 733
$ws.Cells.Item(41, 7).Value = 45
$ws.Cells.Item(41, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=84184'
$ws.Cells.Item(41, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg'

# row 42: 南昌·花绒万兽第二聚
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = '''2024-08-10'
$ws.Cells.Item(42, 3).Value = '''南昌·花绒万兽第二聚'
$ws.Cells.Item(42, 4).Value = '''南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Cells.Item(42, 5).Value = '''2024.08.10 10:00-08.11 17:00'
$ws.Cells.Item(42, 6).Value = 10
$ws.Cells.Item(42, 7).Value = 188
$ws.Cells.Item(42, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87600'
$ws.Cells.Item(42, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/i0Ojsne01718693886054.png'

# row 43: 高安·第二届静卿国风动漫文化展览会
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = '''2024-08-10'
$ws.Cells.Item(43, 3).Value = '''高安·第二届静卿国风动漫文化展览会'
$ws.Cells.Item(43, 4).Value = '''华林中路606号 高安华鼎国际大酒店'
$ws.Cells.Item(43, 5).Value = '''2024.08.10 09:00-08.10 17:00'
$ws.Cells.Item(43, 6).Value = 444
$ws.Cells.Item(43, 7).Value = 45
$ws.Cells.Item(43, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86682'
$ws.Cells.Item(43, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg'

# row 44: 上饶·次元重现夏日嘉年华
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = '''2024-08-15'
$ws.Cells.Item(44, 3).Value = '''上饶·次元重现夏日嘉年华'
$ws.Cells.Item(44, 4).Value = '''普济巷地委大院北侧约90米 四季体育运动馆'
$ws.Cells.Item(44, 5).Value = '''2024.08.15 09:30-08.15 17:30'
$ws.Cells.Item(44, 6).Value = 14
$ws.Cells.Item(44, 7).Value = 48
$ws.Cells.Item(44, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87679'
$ws.Cells.Item(44, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/fxlKV2SL1718784421064.jpeg'

# row 45: 南昌·第四届龙年动漫展——暑假最后的狂欢
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = '''2024-08-24'
$ws.Cells.Item(45, 3).Value = '''南昌·第四届龙年动漫展——暑假最后的狂欢'
$ws.Cells.Item(45, 4).Value = '''南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Cells.Item(45, 5).Value = '''2024.08.24 10:00-08.25 18:00'
$ws.Cells.Item(45, 6).Value = 435
$ws.Cells.Item(45, 7).Value = 45
$ws.Cells.Item(45, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87135'
$ws.Cells.Item(45, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg'


# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Rows.Item(25).Insert() | Out-Null
$ws.Cells.Item(24, 1).Copy() | Out-Null
$ws.Cells.Item(25, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# row 2: 景德镇·BM次元盛典运动番only（取消）
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = '''2024-06-22'
$ws.Cells.Item(2, 3).Value = '''景德镇·BM次元盛典运动番only（取消）'
$ws.Cells.Item(2, 4).Value = '''广场南路金幕影城旁 罗曼园宴会酒店'
$ws.Cells.Item(2, 5).Value = '''2024.06.22 10:00-06.22 17:00'
$ws.Cells.Item(2, 6).Value = 204
$ws.Cells.Item(2, 7).Value = '不可售'
$ws.Cells.Item(2, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85197'
$ws.Cells.Item(2, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202404/Z6eXz0su1714292081978.png'

# row 3: 上饶·BM次元盛典运动番only（取消）
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = '''2024-06-23'
$ws.Cells.Item(3, 3).Value = '''上饶·BM次元盛典运动番only（取消）'
$ws.Cells.Item(3, 4).Value = '''春江北大道时光PARK内 博悦宴会艺术中心'
$ws.Cells.Item(3, 5).Value = '''2024.06.23 10:00-06.23 17:00'
$ws.Cells.Item(3, 6).Value = 279
$ws.Cells.Item(3, 7).Value = '不可售'
$ws.Cells.Item(3, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85201'
$ws.Cells.Item(3, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202404/30dgkbjT1714293499693.png'

# row 4: 赣州·清风霁月·光夜only
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = '''2024-06-23'
$ws.Cells.Item(4, 3).Value = '''赣州·清风霁月·光夜only'
$ws.Cells.Item(4, 4).Value = '''平安大道 麋鹿LiveHouse'
$ws.Cells.Item(4, 5).Value = '''2024.06.23 14:00-06.23 20:00'
$ws.Cells.Item(4, 6).Value = 84
$ws.Cells.Item(4, 7).Value = 158
$ws.Cells.Item(4, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86993'
$ws.Cells.Item(4, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/PklWR8EP1717429316070.jpeg'

# row 5: 南昌·第五人格only
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = '''2024-06-29'
$ws.Cells.Item(5, 3).Value = '''南昌·第五人格only'
$ws.Cells.Item(5, 4).Value = '''高处见美好生活公园 百家喜宴高新店'
$ws.Cells.Item(5, 5).Value = '''2024.06.29 10:00-06.29 17:00'
$ws.Cells.Item(5, 6).Value = 293
$ws.Cells.Item(5, 7).Value = 65
$ws.Cells.Item(5, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87043'
$ws.Cells.Item(5, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202405/zir2PYz81717071721569.jpeg'

# row 6: 萍乡·BM次元盛典运动番only
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = '''2024-06-29'
$ws.Cells.Item(6, 3).Value = '''萍乡·BM次元盛典运动番only'
$ws.Cells.Item(6, 4).Value = '''康庄路3号 萍乡梅园国际大酒店'
$ws.Cells.Item(6, 5).Value = '''2024.06.29 10:00-06.29 17:00'
$ws.Cells.Item(6, 6).Value = 272
$ws.Cells.Item(6, 7).Value = 55
$ws.Cells.Item(6, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85192'
$ws.Cells.Item(6, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png'

# row 7: 南昌·ChinastyleCOSPLAY  
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = '''2024-06-30'
$ws.Cells.Item(7, 3).Value = '''南昌·ChinastyleCOSPLAY  '
$ws.Cells.Item(7, 4).Value = '''真君路888号 南昌华侨城玩美公园'
$ws.Cells.Item(7, 5).Value = '''2024.06.30 09:30-07.02 17:30'
$ws.Cells.Item(7, 6).Value = 116
$ws.Cells.Item(7, 7).Value = 65
$ws.Cells.Item(7, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87045'
$ws.Cells.Item(7, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/wajWy7ID1717149642528.jpeg'

# row 8: 宜春·BM次元盛典运动番only
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = '''2024-06-30'
$ws.Cells.Item(8, 3).Value = '''宜春·BM次元盛典运动番only'
$ws.Cells.Item(8, 4).Value = '''鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)'
$ws.Cells.Item(8, 5).Value = '''2024.06.30 10:00-06.30 17:00'
$ws.Cells.Item(8, 6).Value = 265
$ws.Cells.Item(8, 7).Value = 55
$ws.Cells.Item(8, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=84636'
$ws.Cells.Item(8, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png'

# row 9: 南昌·次元星球动漫游戏展
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = '''2024-07-06'
$ws.Cells.Item(9, 3).Value = '''南昌·次元星球动漫游戏展'
$ws.Cells.Item(9, 4).Value = '''龙蟠街666号融创茂1层 融创茂'
$ws.Cells.Item(9, 5).Value = '''2024.07.06 10:00-07.06 17:00'
$ws.Cells.Item(9, 6).Value = 22
$ws.Cells.Item(9, 7).Value = '不可售'
$ws.Cells.Item(9, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86405'
$ws.Cells.Item(9, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg'

# row 10: 鹰潭·BM次元盛典运动番only
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = '''2024-07-06'
$ws.Cells.Item(10, 3).Value = '''鹰潭·BM次元盛典运动番only'
$ws.Cells.Item(10, 4).Value = '''体育馆东路2号九小隔壁 忆江南•宴会楼'
$ws.Cells.Item(10, 5).Value = '''2024.07.06 10:00-07.06 17:00'
$ws.Cells.Item(10, 6).Value = 52
$ws.Cells.Item(10, 7).Value = 55
$ws.Cells.Item(10, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85997'
$ws.Cells.Item(10, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png'

# row 11: 赣州·BM次元盛典运动番only
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = '''2024-07-07'
$ws.Cells.Item(11, 3).Value = '''赣州·BM次元盛典运动番only'
$ws.Cells.Item(11, 4).Value = '''米瑞金路2口0号上客天下1楼 上客天下.老虔州'
$ws.Cells.Item(11, 5).Value = '''2024.07.07 10:00-07.07 17:00'
$ws.Cells.Item(11, 6).Value = 38
$ws.Cells.Item(11, 7).Value = 55
$ws.Cells.Item(11, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86602'
$ws.Cells.Item(11, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png'

# row 12: 新余·2024第三届MG动漫嘉年华
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = '''2024-07-12'
$ws.Cells.Item(12, 3).Value = '''新余·2024第三届MG动漫嘉年华'
$ws.Cells.Item(12, 4).Value = '''仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅'
$ws.Cells.Item(12, 5).Value = '''2024.07.12 10:00-07.13 17:30'
$ws.Cells.Item(12, 6).Value = 123
$ws.Cells.Item(12, 7).Value = 55
$ws.Cells.Item(12, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86536'
$ws.Cells.Item(12, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg'

# row 13: 南昌·SuperComic动漫游戏博览会
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = '''2024-07-13'
$ws.Cells.Item(13, 3).Value = '''南昌·SuperComic动漫游戏博览会'
$ws.Cells.Item(13, 4).Value = '''怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Cells.Item(13, 5).Value = '''2024.07.13 09:00-07.14 17:00'
$ws.Cells.Item(13, 6).Value = 2547
$ws.Cells.Item(13, 7).Value = 65
$ws.Cells.Item(13, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86992'
$ws.Cells.Item(13, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/wQTAjelJ1717642148929.jpeg'

# row 14: 南昌·SuperComic配音演员刘明月专场见面会
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = '''2024-07-13'
$ws.Cells.Item(14, 3).Value = '''南昌·SuperComic配音演员刘明月专场见面会'
$ws.Cells.Item(14, 4).Value = '''怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Cells.Item(14, 5).Value = '''2024.07.13 09:00-07.13 17:00'
$ws.Cells.Item(14, 6).Value = 88
$ws.Cells.Item(14, 7).Value = 168
$ws.Cells.Item(14, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87570'
$ws.Cells.Item(14, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202406/1D1reIl81718609013880.png'

# row 15: 南昌·THO-梦违赣鄱荟萃·叁~幻想Strawberry~!!
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = '''2024-07-13'
$ws.Cells.Item(15, 3).Value = '''南昌·THO-梦违赣鄱荟萃·叁~幻想Strawberry~!!'
$ws.Cells.Item(15, 4).Value = '''民德路411号 东方豪景花园酒店(民德路店)'
$ws.Cells.Item(15, 5).Value = '''2024.07.13 09:30-07.13 17:30'
$ws.Cells.Item(15, 6).Value = 30
$ws.Cells.Item(15, 7).Value = 65
$ws.Cells.Item(15, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87668'
$ws.Cells.Item(15, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/Bk9cYryT1718360290362.jpeg'

# row 16: 宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = '''2024-07-13'
$ws.Cells.Item(16, 3).Value = '''宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华'
$ws.Cells.Item(16, 4).Value = '''宜春国际商贸城会展中心 宜春国际商贸城会展中心'
$ws.Cells.Item(16, 5).Value = '''2024.07.13 10:00-07.14 17:00'
$ws.Cells.Item(16, 6).Value = 61
$ws.Cells.Item(16, 7).Value = 55
$ws.Cells.Item(16, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86667'
$ws.Cells.Item(16, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg'

# row 17: 江西·东方LiveParty×THO03幻想Strawberry~！！
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = '''2024-07-13'
$ws.Cells.Item(17, 3).Value = '''江西·东方LiveParty×THO03幻想Strawberry~！！'
$ws.Cells.Item(17, 4).Value = '''上海路543号520Park文创公园21号01区域 瓦肆VAS NANCHANG'
$ws.Cells.Item(17, 5).Value = '''2024.07.13 20:30-07.13 23:00'
$ws.Cells.Item(17, 6).Value = 47
$ws.Cells.Item(17, 7).Value = 80
$ws.Cells.Item(17, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87366'
$ws.Cells.Item(17, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/1L3I6Qmg1718292516616.jpeg'

# row 18: 赣州·十万伏特-次元交流会（夏）
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = '''2024-07-13'
$ws.Cells.Item(18, 3).Value = '''赣州·十万伏特-次元交流会（夏）'
$ws.Cells.Item(18, 4).Value = '''梅关大道36-16号 麋鹿星球艺术中心'
$ws.Cells.Item(18, 5).Value = '''2024.07.13 09:30-07.13 17:00'
$ws.Cells.Item(18, 6).Value = 20
$ws.Cells.Item(18, 7).Value = 45
$ws.Cells.Item(18, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87597'
$ws.Cells.Item(18, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/87yQ4Hmf1718681348727.jpeg'

# row 19: 南昌·赛马娘ONLY
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = '''2024-07-14'
$ws.Cells.Item(19, 3).Value = '''南昌·赛马娘ONLY'
$ws.Cells.Item(19, 4).Value = '''洪城路99号 锦都皇冠酒店(八一广场火车站店)'
$ws.Cells.Item(19, 5).Value = '''2024.07.14 09:00-07.14 17:30'
$ws.Cells.Item(19, 6).Value = 42
$ws.Cells.Item(19, 7).Value = 68
$ws.Cells.Item(19, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87367'
$ws.Cells.Item(19, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/wXQuIKtu1718165450704.png'

# row 20: 吉安·COMIC LIFE次元假日05
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = '''2024-07-14'
$ws.Cells.Item(20, 3).Value = '''吉安·COMIC LIFE次元假日05'
$ws.Cells.Item(20, 4).Value = '''东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws.Cells.Item(20, 5).Value = '''2024.07.14 09:00-07.14 18:00'
$ws.Cells.Item(20, 6).Value = 542
$ws.Cells.Item(20, 7).Value = 52.1
$ws.Cells.Item(20, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85924'
$ws.Cells.Item(20, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg'

# row 21: 赣州·第四届赣州半夏动漫展
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = '''2024-07-19'
$ws.Cells.Item(21, 3).Value = '''赣州·第四届赣州半夏动漫展'
$ws.Cells.Item(21, 4).Value = '''105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws.Cells.Item(21, 5).Value = '''2024.07.19 10:00-07.21 17:00'
$ws.Cells.Item(21, 6).Value = 605
$ws.Cells.Item(21, 7).Value = 55
$ws.Cells.Item(21, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86587'
$ws.Cells.Item(21, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg'

# row 22: 南昌·漫拥动漫嘉年华Pro-追光启航
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = '''2024-07-20'
$ws.Cells.Item(22, 3).Value = '''南昌·漫拥动漫嘉年华Pro-追光启航'
$ws.Cells.Item(22, 4).Value = '''小蓝南路420号 洪州体育馆'
$ws.Cells.Item(22, 5).Value = '''2024.07.20 09:00-07.21 17:00'
$ws.Cells.Item(22, 6).Value = 181
$ws.Cells.Item(22, 7).Value = 52.5
$ws.Cells.Item(22, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85796'
$ws.Cells.Item(22, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png'

# row 23: 乐平·CY境界次元动漫夏时庆
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = '''2024-07-21'
$ws.Cells.Item(23, 3).Value = '''乐平·CY境界次元动漫夏时庆'
$ws.Cells.Item(23, 4).Value = '''翥山西路182号 佳佳基大酒店'
$ws.Cells.Item(23, 5).Value = '''2024.07.21 10:00-07.21 17:00'
$ws.Cells.Item(23, 6).Value = 93
$ws.Cells.Item(23, 7).Value = 35
$ws.Cells.Item(23, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86768'
$ws.Cells.Item(23, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png'

# row 24: 九江·SXD动漫嘉年华
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = '''2024-07-21'
$ws.Cells.Item(24, 3).Value = '''九江·SXD动漫嘉年华'
$ws.Cells.Item(24, 4).Value = '''湓浦街道大中路339号 百嘉洲际酒店'
$ws.Cells.Item(24, 5).Value = '''2024.07.21 10:00-07.21 17:30'
$ws.Cells.Item(24, 6).Value = 51
$ws.Cells.Item(24, 7).Value = 45
$ws.Cells.Item(24, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86832'
$ws.Cells.Item(24, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg'

# row 25: 抚州·临次元08·盛夏动漫狂欢节
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = '''2024-07-21'
$ws.Cells.Item(25, 3).Value = '''抚州·临次元08·盛夏动漫狂欢节'
$ws.Cells.Item(25, 4).Value = '''伍塘路1098号 乐课篮球公园'
$ws.Cells.Item(25, 5).Value = '''2024.07.21 10:00-07.21 16:00'
$ws.Cells.Item(25, 6).Value = 11
$ws.Cells.Item(25, 7).Value = 39.9
$ws.Cells.Item(25, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87763'
$ws.Cells.Item(25, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/6qgetbCh1718720523395.jpeg'

# row 26: 萍乡·NL14动漫游戏展·夏日狂想曲
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = '''2024-07-21'
$ws.Cells.Item(26, 3).Value = '''萍乡·NL14动漫游戏展·夏日狂想曲'
$ws.Cells.Item(26, 4).Value = '''公园南路168号(近工行城北分理处) 梅生嘉华酒店'
$ws.Cells.Item(26, 5).Value = '''2024.07.21 10:00-07.21 17:00'
$ws.Cells.Item(26, 6).Value = 55
$ws.Cells.Item(26, 7).Value = 40
$ws.Cells.Item(26, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86658'
$ws.Cells.Item(26, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg'

# row 27: 南昌·萌卡动漫展
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = '''2024-07-26'
$ws.Cells.Item(27, 3).Value = '''南昌·萌卡动漫展'
$ws.Cells.Item(27, 4).Value = '''八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$ws.Cells.Item(27, 5).Value = '''2024.07.26 09:00-07.28 17:00'
$ws.Cells.Item(27, 6).Value = 2141
$ws.Cells.Item(27, 7).Value = 39.9
$ws.Cells.Item(27, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86776'
$ws.Cells.Item(27, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg'

# row 28: 江西·次元星河动漫游戏嘉年华
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = '''2024-07-27'
$ws.Cells.Item(28, 3).Value = '''江西·次元星河动漫游戏嘉年华'
$ws.Cells.Item(28, 4).Value = '''九龙大道1177号 南昌绿地国际博览中心'
$ws.Cells.Item(28, 5).Value = '''2024.07.27 10:00-07.28 17:00'
$ws.Cells.Item(28, 6).Value = 4240
$ws.Cells.Item(28, 7).Value = 69
$ws.Cells.Item(28, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85493'
$ws.Cells.Item(28, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png'

# row 29: 赣州·马娘only
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = '''2024-07-27'
$ws.Cells.Item(29, 3).Value = '''赣州·马娘only'
$ws.Cells.Item(29, 4).Value = '''火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)'
$ws.Cells.Item(29, 5).Value = '''2024.07.27 09:00-07.27 17:00'
$ws.Cells.Item(29, 6).Value = 32
$ws.Cells.Item(29, 7).Value = 60
$ws.Cells.Item(29, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86772'
$ws.Cells.Item(29, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png'

# row 30: 赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = '''2024-07-28'
$ws.Cells.Item(30, 3).Value = '''赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会'
$ws.Cells.Item(30, 4).Value = '''兴国路恒大帝景西门 江西长庚控股有限公司'
$ws.Cells.Item(30, 5).Value = '''2024.07.28 11:00-07.28 17:00'
$ws.Cells.Item(30, 6).Value = 66
$ws.Cells.Item(30, 7).Value = 56
$ws.Cells.Item(30, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85688'
$ws.Cells.Item(30, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png'

# row 31: 宜春·第三十五届静卿国风动漫文化展览会
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = '''2024-07-30'
$ws.Cells.Item(31, 3).Value = '''宜春·第三十五届静卿国风动漫文化展览会'
$ws.Cells.Item(31, 4).Value = '''宜阳大道19号(交通银行旁) 宜春安缦文华酒店'
$ws.Cells.Item(31, 5).Value = '''2024.07.30 09:00-07.30 17:00'
$ws.Cells.Item(31, 6).Value = 467
$ws.Cells.Item(31, 7).Value = 45
$ws.Cells.Item(31, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86684'
$ws.Cells.Item(31, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/45bGPXfQ1716709212619.jpeg'

# row 32: 南昌·幻梦境国际动漫游戏嘉年华1th
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = '''2024-08-03'
$ws.Cells.Item(32, 3).Value = '''南昌·幻梦境国际动漫游戏嘉年华1th'
$ws.Cells.Item(32, 4).Value = '''南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Cells.Item(32, 5).Value = '''2024.08.03 09:00-08.04 17:30'
$ws.Cells.Item(32, 6).Value = 1234
$ws.Cells.Item(32, 7).Value = 64
$ws.Cells.Item(32, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=83980'
$ws.Cells.Item(32, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg'

# row 33: 吉安·COMIC LIFE周年庆典
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = '''2024-08-03'
$ws.Cells.Item(33, 3).Value = '''吉安·COMIC LIFE周年庆典'
$ws.Cells.Item(33, 4).Value = '''东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws.Cells.Item(33, 5).Value = '''2024.08.03 09:30-08.03 18:00'
$ws.Cells.Item(33, 6).Value = 246
$ws.Cells.Item(33, 7).Value = 46.6
$ws.Cells.Item(33, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87164'
$ws.Cells.Item(33, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/NWD9iQ9h1717598526259.jpeg'

# row 34: 景德镇·第十五届瓷都ACG动漫游戏博览会
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = '''2024-08-03'
$ws.Cells.Item(34, 3).Value = '''景德镇·第十五届瓷都ACG动漫游戏博览会'
$ws.Cells.Item(34, 4).Value = '''迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Cells.Item(34, 5).Value = '''2024.08.03 09:00-08.04 17:00'
$ws.Cells.Item(34, 6).Value = 2140
$ws.Cells.Item(34, 7).Value = 55
$ws.Cells.Item(34, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86341'
$ws.Cells.Item(34, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png'

# row 35: 景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = '''2024-08-03'
$ws.Cells.Item(35, 3).Value = '''景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票'
$ws.Cells.Item(35, 4).Value = '''迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Cells.Item(35, 5).Value = '''2024.08.03 08:30-08.03 17:00'
$ws.Cells.Item(35, 6).Value = 568
$ws.Cells.Item(35, 7).Value = '已售罄'
$ws.Cells.Item(35, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=85981'
$ws.Cells.Item(35, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png'

# row 36: 樟树·第二届静卿国风动漫文化展览会
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = '''2024-08-03'
$ws.Cells.Item(36, 3).Value = '''樟树·第二届静卿国风动漫文化展览会'
$ws.Cells.Item(36, 4).Value = '''杏佛路89号 樟树银河国际酒店'
$ws.Cells.Item(36, 5).Value = '''2024.08.03 09:00-08.03 17:00'
$ws.Cells.Item(36, 6).Value = 474
$ws.Cells.Item(36, 7).Value = 45
$ws.Cells.Item(36, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86683'
$ws.Cells.Item(36, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg'

# row 37: 萍乡·AU9夏至国漫展
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = '''2024-08-03'
$ws.Cells.Item(37, 3).Value = '''萍乡·AU9夏至国漫展'
$ws.Cells.Item(37, 4).Value = '''金陵东路18号 萍乡市体育馆'
$ws.Cells.Item(37, 5).Value = '''2024.08.03 10:00-08.03 17:00'
$ws.Cells.Item(37, 6).Value = 69
$ws.Cells.Item(37, 7).Value = 45
$ws.Cells.Item(37, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86453'
$ws.Cells.Item(37, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg'

# row 38: 赣州·第一届环梦动漫游戏嘉年华
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = '''2024-08-03'
$ws.Cells.Item(38, 3).Value = '''赣州·第一届环梦动漫游戏嘉年华'
$ws.Cells.Item(38, 4).Value = '''105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws.Cells.Item(38, 5).Value = '''2024.08.03 09:00-08.05 17:00'
$ws.Cells.Item(38, 6).Value = 25
$ws.Cells.Item(38, 7).Value = 36.6
$ws.Cells.Item(38, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87449'
$ws.Cells.Item(38, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/iC3PuUfR1717740188790.jpeg'

# row 39: 上饶·第十五届IX Group国风嘉年华暨十周年庆典
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = '''2024-08-04'
$ws.Cells.Item(39, 3).Value = '''上饶·第十五届IX Group国风嘉年华暨十周年庆典'
$ws.Cells.Item(39, 4).Value = '''高铁经济试验区凤凰东大道16号7幢 上饶饶商金茂诚悦酒店(上饶高铁站)'
$ws.Cells.Item(39, 5).Value = '''2024.08.04 09:30-08.04 17:30'
$ws.Cells.Item(39, 6).Value = 135
$ws.Cells.Item(39, 7).Value = 60
$ws.Cells.Item(39, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87225'
$ws.Cells.Item(39, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202406/l5fIXZSX1717562269098.jpeg'

# row 40: 九江·第一届异次元动漫嘉年华
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = '''2024-08-04'
$ws.Cells.Item(40, 3).Value = '''九江·第一届异次元动漫嘉年华'
$ws.Cells.Item(40, 4).Value = '''长虹西大道兴城广场99号 九江半岛宾馆'
$ws.Cells.Item(40, 5).Value = '''2024.08.04 08:00-08.04 17:00'
$ws.Cells.Item(40, 6).Value = 296
$ws.Cells.Item(40, 7).Value = 45
$ws.Cells.Item(40, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=84407'
$ws.Cells.Item(40, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202406/65hJjOfJ1717642614493.jpeg'

# row 41: 南昌·第一届异次元动漫嘉年华
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = '''2024-08-06'
$ws.Cells.Item(41, 3).Value = '''南昌·第一届异次元动漫嘉年华'
$ws.Cells.Item(41, 4).Value = '''民德路411号 东方豪景花园酒店(民德路店)'
$ws.Cells.Item(41, 5).Value = '''2024.08.06 08:00-08.06 17:00'
$ws.Cells.Item(41, 6).Value = 442
$ws.Cells.Item(41, 7).Value = 55
$ws.Cells.Item(41, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=84102'
$ws.Cells.Item(41, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg'

# row 42: 赣州·第二届异次元动漫嘉年华
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = '''2024-08-08'
$ws.Cells.Item(42, 3).Value = '''赣州·第二届异次元动漫嘉年华'
$ws.Cells.Item(42, 4).Value = '''金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆'
$ws.Cells.Item(42, 5).Value = '''2024.08.08 08:00-08.08 17:00'
$ws.Cells.Item(42, 6).Value = 733
$ws.Cells.Item(42, 7).Value = 45
$ws.Cells.Item(42, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=84184'
$ws.Cells.Item(42, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg'

# row 43: 南昌·花绒万兽第二聚
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = '''2024-08-10'
$ws.Cells.Item(43, 3).Value = '''南昌·花绒万兽第二聚'
$ws.Cells.Item(43, 4).Value = '''南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Cells.Item(43, 5).Value = '''2024.08.10 10:00-08.11 17:00'
$ws.Cells.Item(43, 6).Value = 10
$ws.Cells.Item(43, 7).Value = 188
$ws.Cells.Item(43, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87600'
$ws.Cells.Item(43, 9).Value = '''//i1.hdslb.com/bfs/openplatform/202406/i0Ojsne01718693886054.png'

# row 44: 高安·第二届静卿国风动漫文化展览会
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = '''2024-08-10'
$ws.Cells.Item(44, 3).Value = '''高安·第二届静卿国风动漫文化展览会'
$ws.Cells.Item(44, 4).Value = '''华林中路606号 高安华鼎国际大酒店'
$ws.Cells.Item(44, 5).Value = '''2024.08.10 09:00-08.10 17:00'
$ws.Cells.Item(44, 6).Value = 444
$ws.Cells.Item(44, 7).Value = 45
$ws.Cells.Item(44, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=86682'
$ws.Cells.Item(44, 9).Value = '''//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg'

# row 45: 上饶·次元重现夏日嘉年华
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = '''2024-08-15'
$ws.Cells.Item(45, 3).Value = '''上饶·次元重现夏日嘉年华'
$ws.Cells.Item(45, 4).Value = '''普济巷地委大院北侧约90米 四季体育运动馆'
$ws.Cells.Item(45, 5).Value = '''2024.08.15 09:30-08.15 17:30'
$ws.Cells.Item(45, 6).Value = 14
$ws.Cells.Item(45, 7).Value = 48
$ws.Cells.Item(45, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87679'
$ws.Cells.Item(45, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/fxlKV2SL1718784421064.jpeg'

# row 46: 南昌·第四届龙年动漫展——暑假最后的狂欢
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = '''2024-08-24'
$ws.Cells.Item(46, 3).Value = '''南昌·第四届龙年动漫展——暑假最后的狂欢'
$ws.Cells.Item(46, 4).Value = '''南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Cells.Item(46, 5).Value = '''2024.08.24 10:00-08.25 18:00'
$ws.Cells.Item(46, 6).Value = 435
$ws.Cells.Item(46, 7).Value = 45
$ws.Cells.Item(46, 8).Value = '''https://show.bilibili.com/platform/detail.html?id=87135'
$ws.Cells.Item(46, 9).Value = '''//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg'

